$d = $word.ActiveDocument

# 1. Ativação date change
$d.Content.Find.Execute("Ativação: 01/01/2020", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ativação: 01/01/2025", 2)

# 2. Portuguese objectives paragraph: "em equipes," -> "em equipes e dentro da empresa,"
$d.Content.Find.Execute("para, em equipes, apresentarem", $true, $false, $false, $false, $false,
                         $true, 1, $false, "para, em equipes e dentro da empresa, apresentarem", 2)

# 3. English objectives paragraph - full replacement
$d.Content.Find.Execute("To lead students to experience in-depth real problems of the industry in order to present the possible solutions in teams, so that they develop transversal skills that are fundamental to their professional life, such as teamwork, project management, pro activity, at the same time in which they consolidate the knowledge acquired during the course.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Guide students to experience real industry problems in more depth so that, in teams and within the company, they present possible solutions, so that they develop fundamental transversal skills for their professional life, such as teamwork, management of projects, pro activity, while consolidating the knowledge acquired during the course.", 2)

# 4. Portuguese "Programa" paragraph - insert text about team work development
$d.Content.Find.Execute("o desenvolvimento das habilidades essenciais para o trabalho em equipes; Inovação Sistemática", $true, $false, $false, $false, $false,
                         $true, 1, $false, "o desenvolvimento das habilidades essenciais para o trabalho em equipes ocorrerá por meio do trabalho em equipes e de reuniões e visitas didáticas realiadas na empresa (mínimo 3);Inovação Sistemática", 2)

# 5. English "Programa" paragraph - full replacement
$d.Content.Find.Execute("Training and work in teams and communication - the development of skills essential to work in teams; Systematic Innovation - development of innovative solutions, systematization and characteristics; Legislation - notions of legislation applied to corporate action; Project Management and Schedule - Methodologies and necessary schematizations with the management elements; Problem Identification - systematization of actions to locate causes; Formulation of the Project - presentation of the managerial aspects necessary for the development of the project, Management Plan, Project Analytical Structure (EAP) etc; Specification of Problems - systematization of problems within the areas of knowledge; Analysis of Available Knowledge, Evaluation and Decision Making; Reporting - formatting within ABNT standards; Presentation of Projects.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Training and work in teams and Communication – the development of essential skills for working in teams will occur through work in teams and meetings and educational visits carried out in the company (minimum 3); Systematic Innovation - development of innovative solutions, systematization and characteristics; Legislation - notions of legislation applied to corporate action; Project Management and Schedule - Methodologies and necessary schematizations with the management elements; Problem Identification - systematization of actions to locate causes; Formulation of the Project - presentation of the managerial aspects necessary for the development of the project, Management Plan, Project Analytical Structure (EAP) etc; Specification of Problems - systematization of problems within the areas of knowledge; Analysis of Available Knowledge, Evaluation and Decision Making; Reporting - formatting within ABNT standards; Presentation of Projects.", 2)

# 6. Evaluation criterion paragraph - full replacement
$d.Content.Find.Execute("Serão feitas duas avaliações por uma banca de professores que assistirão às apresentações, as notas serão as médias das notas dadas pelos professores.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Serão feitas três apresentações pelas equipes realizadas no ambiente físico da empresa parceira, as notas serão compostas pelas avaliações dos tutores da empresa e da Escola.", 2)
